$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) First paragraph: add two trailing spaces to the existing sentence
#    and append a new red run: "(This is a change – Version for branch
#    alternate)" split across three runs, all colored C00000.
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs(1)
$xml1 = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:r><w:t xml:space="preserve">This is a Microsoft word document.  </w:t></w:r>
<w:r><w:rPr><w:color w:val="C00000"/></w:rPr><w:t>(This is a change &#8211; Ve</w:t></w:r>
<w:r><w:rPr><w:color w:val="C00000"/></w:rPr><w:t>rsion for branch alternate</w:t></w:r>
<w:r><w:rPr><w:color w:val="C00000"/></w:rPr><w:t>)</w:t></w:r>
</w:p>
"@
$p1.Range.InsertXML($xml1)

# ---------------------------------------------------------------------
# 2) Fourth paragraph ("Crispian's Day speech ..."): merge the run
#    " Day speech from" with the following stray space run into a
#    single " Day speech from " run (this also relocates the spell-
#    check proofErr markers to wrap only "Shakespear's"), then merge
#    the six runs making up " Henry V" ... "]" into one run reading
#    " Henry V [Source – Wikipedia]".
# ---------------------------------------------------------------------
$dash = [char]8211

$d.Content.Find.Execute(" Day speech from ", $false, $false, $false, $false, `
    $false, $true, 1, $false, " Day speech from ", 2) | Out-Null

$d.Content.Find.Execute(" Henry V [Source " + $dash + " Wikipedia]", $false, $false, $false, $false, `
    $false, $true, 1, $false, " Henry V [Source " + $dash + " Wikipedia]", 2) | Out-Null

# ---------------------------------------------------------------------
# 3) End of document: append two new empty paragraphs before the
#    section break - one carrying the "larger" style plus shading and
#    spacing, and a completely bare one after it.
# ---------------------------------------------------------------------
$xmlEnd = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:pPr>
<w:pStyle w:val="larger"/>
<w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
<w:spacing w:before="0" w:beforeAutospacing="0" w:after="150" w:afterAutospacing="0"/>
</w:pPr>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
"@
$endRange = $d.Range($d.Content.End, $d.Content.End)
$endRange.InsertXML($xmlEnd)
